$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 74-97 with revised prices
$ws.Cells.Item(74, 2).Value = 97.66
$ws.Cells.Item(75, 2).Value = 88.38
$ws.Cells.Item(76, 2).Value = 86.36
$ws.Cells.Item(77, 2).Value = 85.08
$ws.Cells.Item(78, 2).Value = 88.81
$ws.Cells.Item(79, 2).Value = 128.35
$ws.Cells.Item(80, 2).Value = 134.05
$ws.Cells.Item(81, 2).Value = 149.8
$ws.Cells.Item(82, 2).Value = 214.2
$ws.Cells.Item(83, 2).Value = 235.16
$ws.Cells.Item(84, 2).Value = 266.93
$ws.Cells.Item(85, 2).Value = 278.28
$ws.Cells.Item(86, 2).Value = 236.5
$ws.Cells.Item(87, 2).Value = 249.96
$ws.Cells.Item(88, 2).Value = 226.33
$ws.Cells.Item(89, 2).Value = 234.09
$ws.Cells.Item(90, 2).Value = 251.69
$ws.Cells.Item(91, 2).Value = 351.4
$ws.Cells.Item(92, 2).Value = 385.29
$ws.Cells.Item(93, 2).Value = 392.36
$ws.Cells.Item(94, 2).Value = 352.05
$ws.Cells.Item(95, 2).Value = 234.4
$ws.Cells.Item(96, 2).Value = 134.03
$ws.Cells.Item(97, 2).Value = 47.08

# Append new rows 98-169 for 2022-09-26 through 2022-09-28
$ws.Cells.Item(98, 1).Value = "2022-09-26 00:00"
$ws.Cells.Item(98, 2).Value = 137.18
$ws.Cells.Item(99, 1).Value = "2022-09-26 01:00"
$ws.Cells.Item(99, 2).Value = 66.48
$ws.Cells.Item(100, 1).Value = "2022-09-26 02:00"
$ws.Cells.Item(100, 2).Value = 58.56
$ws.Cells.Item(101, 1).Value = "2022-09-26 03:00"
$ws.Cells.Item(101, 2).Value = 61.53
$ws.Cells.Item(102, 1).Value = "2022-09-26 04:00"
$ws.Cells.Item(102, 2).Value = 86.76
$ws.Cells.Item(103, 1).Value = "2022-09-26 05:00"
$ws.Cells.Item(103, 2).Value = 131.65
$ws.Cells.Item(104, 1).Value = "2022-09-26 06:00"
$ws.Cells.Item(104, 2).Value = 64.56
$ws.Cells.Item(105, 1).Value = "2022-09-26 07:00"
$ws.Cells.Item(105, 2).Value = 357.57
$ws.Cells.Item(106, 1).Value = "2022-09-26 08:00"
$ws.Cells.Item(106, 2).Value = 362.08
$ws.Cells.Item(107, 1).Value = "2022-09-26 09:00"
$ws.Cells.Item(107, 2).Value = 319.85
$ws.Cells.Item(108, 1).Value = "2022-09-26 10:00"
$ws.Cells.Item(108, 2).Value = 285.01
$ws.Cells.Item(109, 1).Value = "2022-09-26 11:00"
$ws.Cells.Item(109, 2).Value = 238.63
$ws.Cells.Item(110, 1).Value = "2022-09-26 12:00"
$ws.Cells.Item(110, 2).Value = 178.83
$ws.Cells.Item(111, 1).Value = "2022-09-26 13:00"
$ws.Cells.Item(111, 2).Value = 186.27
$ws.Cells.Item(112, 1).Value = "2022-09-26 14:00"
$ws.Cells.Item(112, 2).Value = 102.67
$ws.Cells.Item(113, 1).Value = "2022-09-26 15:00"
$ws.Cells.Item(113, 2).Value = 67.8
$ws.Cells.Item(114, 1).Value = "2022-09-26 16:00"
$ws.Cells.Item(114, 2).Value = 56.51
$ws.Cells.Item(115, 1).Value = "2022-09-26 17:00"
$ws.Cells.Item(115, 2).Value = 61.76
$ws.Cells.Item(116, 1).Value = "2022-09-26 18:00"
$ws.Cells.Item(116, 2).Value = 70.08
$ws.Cells.Item(117, 1).Value = "2022-09-26 19:00"
$ws.Cells.Item(117, 2).Value = 69.75
$ws.Cells.Item(118, 1).Value = "2022-09-26 20:00"
$ws.Cells.Item(118, 2).Value = 58.68
$ws.Cells.Item(119, 1).Value = "2022-09-26 21:00"
$ws.Cells.Item(119, 2).Value = 46.51
$ws.Cells.Item(120, 1).Value = "2022-09-26 22:00"
$ws.Cells.Item(120, 2).Value = 34.66
$ws.Cells.Item(121, 1).Value = "2022-09-26 23:00"
$ws.Cells.Item(121, 2).Value = 20.6
$ws.Cells.Item(122, 1).Value = "2022-09-27 00:00"
$ws.Cells.Item(122, 2).Value = 15.22
$ws.Cells.Item(123, 1).Value = "2022-09-27 01:00"
$ws.Cells.Item(123, 2).Value = 14.21
$ws.Cells.Item(124, 1).Value = "2022-09-27 02:00"
$ws.Cells.Item(124, 2).Value = 13.84
$ws.Cells.Item(125, 1).Value = "2022-09-27 03:00"
$ws.Cells.Item(125, 2).Value = 14.18
$ws.Cells.Item(126, 1).Value = "2022-09-27 04:00"
$ws.Cells.Item(126, 2).Value = 15.65
$ws.Cells.Item(127, 1).Value = "2022-09-27 05:00"
$ws.Cells.Item(127, 2).Value = 17.71
$ws.Cells.Item(128, 1).Value = "2022-09-27 06:00"
$ws.Cells.Item(128, 2).Value = 27.02
$ws.Cells.Item(129, 1).Value = "2022-09-27 07:00"
$ws.Cells.Item(129, 2).Value = 64.93
$ws.Cells.Item(130, 1).Value = "2022-09-27 08:00"
$ws.Cells.Item(130, 2).Value = 98.54
$ws.Cells.Item(131, 1).Value = "2022-09-27 09:00"
$ws.Cells.Item(131, 2).Value = 76.45
$ws.Cells.Item(132, 1).Value = "2022-09-27 10:00"
$ws.Cells.Item(132, 2).Value = 75.73
$ws.Cells.Item(133, 1).Value = "2022-09-27 11:00"
$ws.Cells.Item(133, 2).Value = 75.15
$ws.Cells.Item(134, 1).Value = "2022-09-27 12:00"
$ws.Cells.Item(134, 2).Value = 67.83
$ws.Cells.Item(135, 1).Value = "2022-09-27 13:00"
$ws.Cells.Item(135, 2).Value = 66.52
$ws.Cells.Item(136, 1).Value = "2022-09-27 14:00"
$ws.Cells.Item(136, 2).Value = 63.92
$ws.Cells.Item(137, 1).Value = "2022-09-27 15:00"
$ws.Cells.Item(137, 2).Value = 62.76
$ws.Cells.Item(138, 1).Value = "2022-09-27 16:00"
$ws.Cells.Item(138, 2).Value = 63.45
$ws.Cells.Item(139, 1).Value = "2022-09-27 17:00"
$ws.Cells.Item(139, 2).Value = 67.82
$ws.Cells.Item(140, 1).Value = "2022-09-27 18:00"
$ws.Cells.Item(140, 2).Value = 72.5
$ws.Cells.Item(141, 1).Value = "2022-09-27 19:00"
$ws.Cells.Item(141, 2).Value = 71.22
$ws.Cells.Item(142, 1).Value = "2022-09-27 20:00"
$ws.Cells.Item(142, 2).Value = 66
$ws.Cells.Item(143, 1).Value = "2022-09-27 21:00"
$ws.Cells.Item(143, 2).Value = 57.24
$ws.Cells.Item(144, 1).Value = "2022-09-27 22:00"
$ws.Cells.Item(144, 2).Value = 50.39
$ws.Cells.Item(145, 1).Value = "2022-09-27 23:00"
$ws.Cells.Item(145, 2).Value = 37.96
$ws.Cells.Item(146, 1).Value = "2022-09-28 00:00"
$ws.Cells.Item(146, 2).Value = 49.5
$ws.Cells.Item(147, 1).Value = "2022-09-28 01:00"
$ws.Cells.Item(147, 2).Value = 49.06
$ws.Cells.Item(148, 1).Value = "2022-09-28 02:00"
$ws.Cells.Item(148, 2).Value = 49.76
$ws.Cells.Item(149, 1).Value = "2022-09-28 03:00"
$ws.Cells.Item(149, 2).Value = 50.72
$ws.Cells.Item(150, 1).Value = "2022-09-28 04:00"
$ws.Cells.Item(150, 2).Value = 51.09
$ws.Cells.Item(151, 1).Value = "2022-09-28 05:00"
$ws.Cells.Item(151, 2).Value = 54.99
$ws.Cells.Item(152, 1).Value = "2022-09-28 06:00"
$ws.Cells.Item(152, 2).Value = 61.35
$ws.Cells.Item(153, 1).Value = "2022-09-28 07:00"
$ws.Cells.Item(153, 2).Value = 74.97
$ws.Cells.Item(154, 1).Value = "2022-09-28 08:00"
$ws.Cells.Item(154, 2).Value = 82.83
$ws.Cells.Item(155, 1).Value = "2022-09-28 09:00"
$ws.Cells.Item(155, 2).Value = 83.9
$ws.Cells.Item(156, 1).Value = "2022-09-28 10:00"
$ws.Cells.Item(156, 2).Value = 83.08
$ws.Cells.Item(157, 1).Value = "2022-09-28 11:00"
$ws.Cells.Item(157, 2).Value = 82.92
$ws.Cells.Item(158, 1).Value = "2022-09-28 12:00"
$ws.Cells.Item(158, 2).Value = 82.65
$ws.Cells.Item(159, 1).Value = "2022-09-28 13:00"
$ws.Cells.Item(159, 2).Value = 83.02
$ws.Cells.Item(160, 1).Value = "2022-09-28 14:00"
$ws.Cells.Item(160, 2).Value = 83.43
$ws.Cells.Item(161, 1).Value = "2022-09-28 15:00"
$ws.Cells.Item(161, 2).Value = 83.5
$ws.Cells.Item(162, 1).Value = "2022-09-28 16:00"
$ws.Cells.Item(162, 2).Value = 78.44
$ws.Cells.Item(163, 1).Value = "2022-09-28 17:00"
$ws.Cells.Item(163, 2).Value = 82.46
$ws.Cells.Item(164, 1).Value = "2022-09-28 18:00"
$ws.Cells.Item(164, 2).Value = 82.24
$ws.Cells.Item(165, 1).Value = "2022-09-28 19:00"
$ws.Cells.Item(165, 2).Value = 82.13
$ws.Cells.Item(166, 1).Value = "2022-09-28 20:00"
$ws.Cells.Item(166, 2).Value = 83.38
$ws.Cells.Item(167, 1).Value = "2022-09-28 21:00"
$ws.Cells.Item(167, 2).Value = 76.3
$ws.Cells.Item(168, 1).Value = "2022-09-28 22:00"
$ws.Cells.Item(168, 2).Value = 76.03
$ws.Cells.Item(169, 1).Value = "2022-09-28 23:00"
$ws.Cells.Item(169, 2).Value = 65.36

# Resize the table to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B170"))

Write-Output "done"
